$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 0.9

$ws1.Range("D3").Value = 27
$ws1.Range("H3").Value = 11.64
$ws1.Range("L3").Value = 0.98

$ws1.Range("D4").Value = 24
$ws1.Range("H4").Value = 12.13
$ws1.Range("L4").Value = 0.96

$ws1.Range("D5").Value = 29
$ws1.Range("H5").Value = 9.17
$ws1.Range("L5").Value = 1.15

$ws1.Range("D6").Value = 29
$ws1.Range("H6").Value = 8.25
$ws1.Range("L6").Value = 0.91

$ws1.Range("H7").Value = 7.31
$ws1.Range("L7").Value = 1.01

$ws1.Range("D8").Value = 29
$ws1.Range("H8").Value = 6.26
$ws1.Range("L8").Value = 1.11

$ws1.Range("D9").Value = 29
$ws1.Range("H9").Value = 5.26
$ws1.Range("L9").Value = 1.02

$ws1.Range("H10").Value = 4.33
$ws1.Range("L10").Value = 0.99

$ws1.Range("D11").Value = 23
$ws1.Range("H11").Value = 4.04
$ws1.Range("L11").Value = 0.98

$ws1.Range("D12").Value = 26
$ws1.Range("H12").Value = 2.67
$ws1.Range("L12").Value = 0.8

$ws1.Range("H13").Value = 1.57
$ws1.Range("L13").Value = 1.03

$ws1.Range("H14").Value = 0.57
$ws1.Range("I14").Value = "Low"
$ws1.Range("L14").Value = 1.02

$ws1.Range("L15").Value = 0.97

$ws1.Range("L16").Value = 0.84

$ws1.Range("L17").Value = 1.05

# --- Sheet "Summary" ---
# These cells store numeric-looking values as TEXT in the workbook, so a
# leading apostrophe is used to force text entry (matching the original
# inlineStr / string cell type) instead of letting them be auto-converted
# to numbers.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'447"
$ws2.Range("B10").Value = "'226"
$ws2.Range("B11").Value = "'110"
$ws2.Range("B14").Value = "'24"
